# LOQ4265.xlsx update
#
# The two rows that only carried the "responsible professors" values (rows
# 13 & 14, which had no label in column A) are removed, so every row below
# them shifts up by two. A handful of the remaining value cells (columns B
# and C) are then updated to their new text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the former rows 13 and 14 — everything below shifts up two rows.
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(13).Delete()

# "Objetivos:" value now holds the professor's name instead of the long text.
$ws.Range("B10").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C10").Value = "5840560 - Marco Antonio Carvalho Pereira"

# "Programa resumido:" value becomes "01/01/2021". Stage it in a scratch
# formula cell first and paste only the resulting value in, so Excel's
# automatic date detection doesn't turn the literal text into a date
# serial number (and pull in a new date number-format in the process).
$ws.Range("Z1").Formula = "=""01/01/2021"""
$ws.Range("Z1").Copy()
$ws.Range("B13").PasteSpecial(-4163)
$ws.Range("C13").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("Z1").ClearContents()

# "Programa:" value.
$ws.Range("B15").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C15").Value = "5840560 - Marco Antonio Carvalho Pereira"

# "Método:" value.
$ws.Range("B18").Value = "1285870 - Marcos Villela Barcza"
$ws.Range("C18").Value = "1285870 - Marcos Villela Barcza"

# "Critério:" value.
$ws.Range("B19").Value = "Reuniões periódicas com o orientador e realização do trabalho de conclusão de curso conforme orientação e apresentação de uma monografia final."
$ws.Range("C19").Value = "Reuniões periódicas com o orientador e realização do trabalho de conclusão de curso conforme orientação e apresentação de uma monografia final."

# "Norma de recuperação:" value.
$ws.Range("B20").Value = "Uma única prova perante uma banca com 3 examinadores. A nota da disciplina será decidida pelos docentes da banca."
$ws.Range("C20").Value = "Uma única prova perante uma banca com 3 examinadores. A nota da disciplina será decidida pelos docentes da banca."

# "Bibliografia:" value.
$ws.Range("B21").Value = "Reapresentação do trabalho modificado para nova avaliação."
$ws.Range("C21").Value = "Reapresentação do trabalho modificado para nova avaliação."
